$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.317.30'
$ws.Range("E2").Value = '  -4.30%  '
$ws.Range("E3").Value = '  -5.43%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.89%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.514'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.02%  '
$ws.Range("D9").Value = '2.501.02'
$ws.Range("E9").Value = '  -5.44%  '
$ws.Range("E10").Value = '  -9.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.167'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.31%  '
$ws.Range("E12").Value = '  -4.10%  '
$ws.Range("E13").Value = '  -2.35%  '
$ws.Range("D14").Value = '2.962.79'
$ws.Range("E14").Value = '  -5.35%  '
$ws.Range("D15").Value = '69.192.76'
$ws.Range("E15").Value = '  -4.25%  '
$ws.Range("E16").Value = '  -7.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.65'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.91%  '
$ws.Range("D18").Value = '2.504.93'
$ws.Range("E18").Value = '  -5.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.74'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '345.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -6.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.76%  '
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.93'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.26%  '
$ws.Range("E25").Value = '  -3.98%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.94'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -7.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = '2.633.62'
$ws.Range("E28").Value = '  -5.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("D30").Value = '0.0₃0891'
$ws.Range("E30").Value = '  -6.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.83'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '461.96'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.99%  '
$ws.Range("E33").Value = '  -2.26%  '
$ws.Range("E34").Value = '  -3.49%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '152.93'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.90%  '
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.32'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.01%  '
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.71'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.314'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.94%  '
$ws.Range("E43").Value = '  -8.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.15'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -15.07%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '38.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.62%  '
$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.28'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -11.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.96%  '
$ws.Range("E48").Value = '  -4.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.90%  '
$ws.Range("E50").Value = '  -5.51%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0729'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.68%  '
